$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text (not numbers) in this sheet. Where the new
# price string would otherwise be auto-parsed as a number by Excel (e.g. "1.00",
# "35.40" -> 35.4), force the cell to Text format first so the literal digits are
# preserved exactly as scraped.

$ws.Range("D2").Value = '41.089.23'
$ws.Range("E2").Value = '  -2.31%  '
$ws.Range("D3").Value = '2.174.02'
$ws.Range("E3").Value = '  -2.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.11'
$ws.Range("E5").Value = '  -1.11%  '
$ws.Range("E6").Value = '  -1.94%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '65.24'
$ws.Range("E7").Value = '  -9.29%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.563'
$ws.Range("E9").Value = '  -5.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.76'
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0925'
$ws.Range("E11").Value = '  -5.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '35.40'
$ws.Range("E12").Value = '  -14.71%  '
$ws.Range("E13").Value = '  -2.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.85'
$ws.Range("E14").Value = '  -5.80%  '
$ws.Range("D15").Value = '2.494.07'
$ws.Range("E15").Value = '  -2.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.26'
$ws.Range("E16").Value = '  -5.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.847'
$ws.Range("E17").Value = '  -2.67%  '
$ws.Range("D18").Value = '2.178.49'
$ws.Range("E18").Value = '  -2.33%  '
$ws.Range("D19").Value = '41.078.87'
$ws.Range("E19").Value = '  -2.04%  '
$ws.Range("D20").Value = '0.0₃0936'
$ws.Range("E20").Value = '  -4.37%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.07'
$ws.Range("E21").Value = '  -2.82%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.42'
$ws.Range("E22").Value = '  -2.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '229.34'
$ws.Range("E23").Value = '  -3.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.04'
$ws.Range("E24").Value = '  -5.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.81'
$ws.Range("E25").Value = '  -9.81%  '
$ws.Range("E26").Value = '  +0.06%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.28'
$ws.Range("E27").Value = '  +3.80%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.40'
$ws.Range("E28").Value = '  -6.07%  '
$ws.Range("B29").Value = 'LEO'
$ws.Range("C29").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.72'
$ws.Range("E29").Value = '  -5.62%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.10'
$ws.Range("E30").Value = '  -4.82%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.89'
$ws.Range("E31").Value = '  -2.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.12'
$ws.Range("E32").Value = '  -3.72%  '
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0737'
$ws.Range("E35").Value = '  +0.83%  '
$ws.Range("E36").Value = '  -4.16%  '
$ws.Range("E37").Value = '  -4.92%  '
$ws.Range("E38").Value = '  +1.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.24'
$ws.Range("E39").Value = '  -7.57%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0303'
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("E41").Value = '  -5.63%  '
$ws.Range("E42").Value = '  -9.28%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.95'
$ws.Range("E43").Value = '  +0.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '60.11'
$ws.Range("E44").Value = '  -12.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '11.13'
$ws.Range("E45").Value = '  -6.42%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.50'
$ws.Range("E46").Value = '  -3.96%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.189'
$ws.Range("E47").Value = '  -9.26%  '
$ws.Range("B48").Value = 'BinanceUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0989'
$ws.Range("E49").Value = '  -3.65%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.14'
$ws.Range("E50").Value = '  -1.67%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.14'
$ws.Range("E51").Value = '  -4.68%  '

Write-Host "Updated cryptos list"